# edit.ps1 - apply "edited seq diagram for visitweb" changes to
# VisitWebSequenceDiagram.pptx (slide 1), matching the authoritative
# OOXML diff as closely as the PowerPoint COM object model allows.
#
# Summary of content changes applied here (see diff for full context):
#   1) Straight Arrow Connector 22 (id 23): line color 00B050 -> 0070C0
#   2) Straight Connector 67 (id 68): un-flip + reposition/resize
#   3) Rectangle 62 (id 71): reposition (Left/Top)
#   4) Straight Arrow Connector 71 (id 72): line color 7030A0 -> 00B050
#   5) Straight Connector 72 (id 73): un-flip + reposition/resize
#   6) Rectangle 74 (id 75) "loadPage(Restaurant.weblink)" textbox:
#        - give it a no-fill line (<a:ln><a:noFill/></a:ln>)
#        - recolor all runs from themed accent4/lumMod75% to solid 00B050
#   (Date placeholder "today" fields on the layouts/master/notesMaster are
#   auto-generated datetimeFigureOut fields; they are intentionally left
#   alone here because writing to them through this object model collapses
#   the <a:fld> into a plain run, which would be a worse structural match
#   than leaving the cached text as-is.)

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shapes = $s.Shapes

# --- 1) id=23 "Straight Arrow Connector 22": ln color 00B050 -> 0070C0 ---
$shp23 = Get-ShapeById $shapes 23
$shp23.Line.ForeColor.RGB = 12611584   # 0x0070C0 in BGR packing

# --- 2) id=68 "Straight Connector 67": remove flipH, new off/ext ---
$shp68 = Get-ShapeById $shapes 68
$shp68.HorizontalFlip = $false
$shp68.Left = 569.96732    # 7238585 EMU
$shp68.Top = 249.2645      # 3165659 EMU (unchanged)
$shp68.Width = 0.03268     # 415 EMU
$shp68.Height = 125.8841   # 1598728 EMU

# --- 3) id=71 "Rectangle 62": reposition only (ext unchanged) ---
$shp71 = Get-ShapeById $shapes 71
$shp71.Left = 660.60648    # 8389702 EMU
$shp71.Top = 226.18583     # 2872560 EMU

# --- 4) id=72 "Straight Arrow Connector 71": ln color 7030A0 -> 00B050 ---
$shp72 = Get-ShapeById $shapes 72
$shp72.Line.ForeColor.RGB = 5287936    # 0x00B050 in BGR packing

# --- 5) id=73 "Straight Connector 72": remove flipH, new off/ext ---
$shp73 = Get-ShapeById $shapes 73
$shp73.HorizontalFlip = $false
$shp73.Left = 717.18314    # 9108226 EMU
$shp73.Top = 254.18693     # 3228174 EMU
$shp73.Width = 1.41111     # 17921 EMU
$shp73.Height = 121.06607  # 1537539 EMU

# --- 6) id=75 "Rectangle 74" loadPage(...) textbox ---
$shp75 = Get-ShapeById $shapes 75
$shp75.Line.Visible = $false   # emits <a:ln><a:noFill/></a:ln>

$tr75 = $shp75.TextFrame.TextRange
# "loadPage" "(" "Restaurant.weblink" ")" -> 4 runs, 8+1+18+1 = 28 chars
$tr75.Characters(1, 8).Font.Color.RGB = 5287936    # loadPage
$tr75.Characters(9, 1).Font.Color.RGB = 5287936    # (
$tr75.Characters(10, 18).Font.Color.RGB = 5287936  # Restaurant.weblink
$tr75.Characters(28, 1).Font.Color.RGB = 5287936   # )
